# Updates TPM-derived values for Efnb1-Ephb4 LR-pair sheet (commit: "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.546140333333334
$ws.Range("H2").Value = 28.638421
$ws.Range("I2").Value = 0.587227294878132
$ws.Range("J2").Value = 0.587227294878132
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 47.991936
$ws.Range("N2").Value = 143.975808
$ws.Range("O2").Value = 0.7605119179168339
$ws.Range("P2").Value = 0.7605119179168338
$ws.Range("Q2").Value = 458.1377559243521
$ws.Range("R2").Value = 4123.239803319168
$ws.Range("S2").Value = 0.4465933562808823
$ws.Range("T2").Value = 0.4465933562808823
$ws.Range("G3").Value = 9.546140333333334
$ws.Range("H3").Value = 28.638421
$ws.Range("I3").Value = 0.587227294878132
$ws.Range("J3").Value = 0.587227294878132
$ws.Range("O3").Value = 0.1317597634642934
$ws.Range("P3").Value = 0.1317597634642934
$ws.Range("Q3").Value = 79.37301300945022
$ws.Range("R3").Value = 714.3571170850521
$ws.Range("S3").Value = 0.07737292947291953
$ws.Range("T3").Value = 0.07737292947291953
$ws.Range("G4").Value = 9.546140333333334
$ws.Range("H4").Value = 28.638421
$ws.Range("I4").Value = 0.587227294878132
$ws.Range("J4").Value = 0.587227294878132
$ws.Range("M4").Value = 6.744108333333334
$ws.Range("N4").Value = 20.232325
$ws.Range("O4").Value = 0.1068715953284784
$ws.Range("P4").Value = 0.1068715953284784
$ws.Range("Q4").Value = 64.38020457320279
$ws.Range("R4").Value = 579.4218411588251
$ws.Range("S4").Value = 0.06275791782405279
$ws.Range("T4").Value = 0.06275791782405279
$ws.Range("G5").Value = 9.546140333333334
$ws.Range("H5").Value = 28.638421
$ws.Range("I5").Value = 0.587227294878132
$ws.Range("J5").Value = 0.587227294878132
$ws.Range("M5").Value = 0.05406333333333333
$ws.Range("N5").Value = 0.16219
$ws.Range("O5").Value = 0.0008567232903942534
$ws.Range("P5").Value = 0.0008567232903942534
$ws.Range("Q5").Value = 0.5160961668877778
$ws.Range("R5").Value = 4.64486550199
$ws.Range("S5").Value = 0.0005030913002773097
$ws.Range("T5").Value = 0.0005030913002773097
$ws.Range("I6").Value = 0.2496684258894083
$ws.Range("J6").Value = 0.2496684258894083
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 47.991936
$ws.Range("N6").Value = 143.975808
$ws.Range("O6").Value = 0.7605119179168339
$ws.Range("P6").Value = 0.7605119179168338
$ws.Range("Q6").Value = 194.784086774912
$ws.Range("R6").Value = 1753.056780974208
$ws.Range("S6").Value = 0.1898758134164308
$ws.Range("T6").Value = 0.1898758134164308
$ws.Range("I7").Value = 0.2496684258894083
$ws.Range("J7").Value = 0.2496684258894083
$ws.Range("O7").Value = 0.1317597634642934
$ws.Range("P7").Value = 0.1317597634642934
$ws.Range("S7").Value = 0.03289625273969089
$ws.Range("T7").Value = 0.03289625273969089
$ws.Range("I8").Value = 0.2496684258894083
$ws.Range("J8").Value = 0.2496684258894083
$ws.Range("M8").Value = 6.744108333333334
$ws.Range("N8").Value = 20.232325
$ws.Range("O8").Value = 0.1068715953284784
$ws.Range("P8").Value = 0.1068715953284784
$ws.Range("Q8").Value = 27.37220233873056
$ws.Range("R8").Value = 246.3498210485751
$ws.Range("S8").Value = 0.02668246297795105
$ws.Range("T8").Value = 0.02668246297795105
$ws.Range("I9").Value = 0.2496684258894083
$ws.Range("J9").Value = 0.2496684258894083
$ws.Range("M9").Value = 0.05406333333333333
$ws.Range("N9").Value = 0.16219
$ws.Range("O9").Value = 0.0008567232903942534
$ws.Range("P9").Value = 0.0008567232903942534
$ws.Range("Q9").Value = 0.2194259679655556
$ws.Range("R9").Value = 1.97483371169
$ws.Range("S9").Value = 0.0002138967553355277
$ws.Range("T9").Value = 0.0002138967553355277
$ws.Range("G10").Value = 2.210442
$ws.Range("H10").Value = 6.631326
$ws.Range("I10").Value = 0.1359745227725727
$ws.Range("J10").Value = 0.1359745227725727
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 47.991936
$ws.Range("N10").Value = 143.975808
$ws.Range("O10").Value = 0.7605119179168339
$ws.Range("P10").Value = 0.7605119179168338
$ws.Range("Q10").Value = 106.083390995712
$ws.Range("R10").Value = 954.750518961408
$ws.Range("S10").Value = 0.1034102451015954
$ws.Range("T10").Value = 0.1034102451015954
$ws.Range("G11").Value = 2.210442
$ws.Range("H11").Value = 6.631326
$ws.Range("I11").Value = 0.1359745227725727
$ws.Range("J11").Value = 0.1359745227725727
$ws.Range("O11").Value = 0.1317597634642934
$ws.Range("P11").Value = 0.1317597634642934
$ws.Range("Q11").Value = 18.379097257768
$ws.Range("R11").Value = 165.411875319912
$ws.Range("S11").Value = 0.01791597095768435
$ws.Range("T11").Value = 0.01791597095768435
$ws.Range("G12").Value = 2.210442
$ws.Range("H12").Value = 6.631326
$ws.Range("I12").Value = 0.1359745227725727
$ws.Range("J12").Value = 0.1359745227725727
$ws.Range("M12").Value = 6.744108333333334
$ws.Range("N12").Value = 20.232325
$ws.Range("O12").Value = 0.1068715953284784
$ws.Range("P12").Value = 0.1068715953284784
$ws.Range("Q12").Value = 14.90746031255
$ws.Range("R12").Value = 134.16714281295
$ws.Range("S12").Value = 0.01453181417273336
$ws.Range("T12").Value = 0.01453181417273336
$ws.Range("G13").Value = 2.210442
$ws.Range("H13").Value = 6.631326
$ws.Range("I13").Value = 0.1359745227725727
$ws.Range("J13").Value = 0.1359745227725727
$ws.Range("M13").Value = 0.05406333333333333
$ws.Range("N13").Value = 0.16219
$ws.Range("O13").Value = 0.0008567232903942534
$ws.Range("P13").Value = 0.0008567232903942534
$ws.Range("Q13").Value = 0.11950386266
$ws.Range("R13").Value = 1.07553476394
$ws.Range("S13").Value = 0.0001164925405595068
$ws.Range("T13").Value = 0.0001164925405595068
$ws.Range("G14").Value = 0.4410293333333333
$ws.Range("H14").Value = 1.323088
$ws.Range("I14").Value = 0.02712975645988715
$ws.Range("J14").Value = 0.02712975645988715
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 47.991936
$ws.Range("N14").Value = 143.975808
$ws.Range("O14").Value = 0.7605119179168339
$ws.Range("P14").Value = 0.7605119179168338
$ws.Range("Q14").Value = 21.165851539456
$ws.Range("R14").Value = 190.492663855104
$ws.Range("S14").Value = 0.02063250311792539
$ws.Range("T14").Value = 0.02063250311792539
$ws.Range("G15").Value = 0.4410293333333333
$ws.Range("H15").Value = 1.323088
$ws.Range("I15").Value = 0.02712975645988715
$ws.Range("J15").Value = 0.02712975645988715
$ws.Range("O15").Value = 0.1317597634642934
$ws.Range("P15").Value = 0.1317597634642934
$ws.Range("Q15").Value = 3.667013661006222
$ws.Range("R15").Value = 33.003122949056
$ws.Range("S15").Value = 0.003574610293998616
$ws.Range("T15").Value = 0.003574610293998616
$ws.Range("G16").Value = 0.4410293333333333
$ws.Range("H16").Value = 1.323088
$ws.Range("I16").Value = 0.02712975645988715
$ws.Range("J16").Value = 0.02712975645988715
$ws.Range("M16").Value = 6.744108333333334
$ws.Range("N16").Value = 20.232325
$ws.Range("O16").Value = 0.1068715953284784
$ws.Range("P16").Value = 0.1068715953284784
$ws.Range("Q16").Value = 2.974349602177778
$ws.Range("R16").Value = 26.76914641960001
$ws.Range("S16").Value = 0.002899400353741233
$ws.Range("T16").Value = 0.002899400353741233
$ws.Range("G17").Value = 0.4410293333333333
$ws.Range("H17").Value = 1.323088
$ws.Range("I17").Value = 0.02712975645988715
$ws.Range("J17").Value = 0.02712975645988715
$ws.Range("M17").Value = 0.05406333333333333
$ws.Range("N17").Value = 0.16219
$ws.Range("O17").Value = 0.0008567232903942534
$ws.Range("P17").Value = 0.0008567232903942534
$ws.Range("Q17").Value = 0.02384351585777778
$ws.Range("R17").Value = 0.21459164272
$ws.Range("S17").Value = 0.00002324269422190927
$ws.Range("T17").Value = 0.00002324269422190927

Write-Host "Updated $($ws.Name) with new TPM values"
